$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, matching the style of the existing
# header cells (e.g. G1: bold font, thin border, centered alignment)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new column's data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
